$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet ---
$ws.Name = "Sheet1"

# --- Header row (row 1): new, wider set of column headers ---
$headers = @(
    "Date",
    "Model Name",
    "Exact Precision (Micro Avg)",
    "Exact Recall (Micro Avg)",
    "Exact F1 Score (Micro Avg)",
    "Exact Precision (Macro Avg)",
    "Exact Recall (Macro Avg)",
    "Exact F1 Score (Macro Avg)",
    "Exact Precision (Weighted Avg)",
    "Exact Recall (Weighted Avg)",
    "Exact F1 Score (Weighted Avg)",
    "Partial Precision",
    "Partial Recall",
    "Partial F1 Score",
    "Partial TP",
    "Partial FP",
    "Partial FN",
    "Support",
    "Accuracy",
    "Result Link",
    "Stats Link",
    "No of GPU Used",
    "Power Consumption"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Style the header row (A1:W1) in one shot: bold font, thin box border,
# centered horizontally and top-aligned vertically.
$headerRange = $ws.Range("A1:W1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# --- Data row (row 2) ---
# A2 holds a date-looking value that must stay plain text, so force it
# through the quote-prefix (leading apostrophe) path instead of letting
# the host auto-convert it to a date serial.
$ws.Range("A2").Value = "'09/10/2025"

$ws.Cells.Item(2, 2).Value = "Llama-3.3-70B-Instruct"
$ws.Cells.Item(2, 3).Value = 0.3333333333333333
$ws.Cells.Item(2, 4).Value = 0.2519083969465649
$ws.Cells.Item(2, 5).Value = 0.2869565217391304
$ws.Cells.Item(2, 6).Value = 0.172611531986532
$ws.Cells.Item(2, 7).Value = 0.09623918074622298
$ws.Cells.Item(2, 8).Value = 0.1230062850262559
$ws.Cells.Item(2, 9).Value = 0.4168433700876449
$ws.Cells.Item(2, 10).Value = 0.2519083969465649
$ws.Cells.Item(2, 11).Value = 0.3121387725415884
$ws.Cells.Item(2, 12).Value = 0.404040404040404
$ws.Cells.Item(2, 13).Value = 0.3065134099616858
$ws.Cells.Item(2, 14).Value = 0.3485838779956427
$ws.Cells.Item(2, 15).Value = 80
$ws.Cells.Item(2, 16).Value = 118
$ws.Cells.Item(2, 17).Value = 181
$ws.Cells.Item(2, 18).Value = 262
$ws.Cells.Item(2, 19).Value = 0.9434835898134277
$ws.Cells.Item(2, 20).Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/ner_evaluation_results_Llama-3.3-70B-Instruct_4_shot.txt"
$ws.Cells.Item(2, 21).Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/Stats/ner_evaluation_stats_Llama-3.3-70B-Instruct_4_shot.txt"
$ws.Cells.Item(2, 22).Value = "4 MLGPU"
$ws.Cells.Item(2, 23).Value = "0.206 kWh"
$ws.Cells.Item(2, 24).Value = 6711

Write-Host "edit applied"
